# Apply the changes described by the diff: update soft-skills/project scores
# and the reasoning chain entries (rows 7-9 now belong to id=3, with new
# project names/scores), and remove the now-redundant rows 10-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-6 (id = 1 group) with the new project_score / final_score values ---
$ws.Range("C2").Value = 112.64
$ws.Range("C3").Value = 102.85
$ws.Range("C4").Value = 83.73999999999999
$ws.Range("E4").Value = 71.18000000000001
$ws.Range("C5").Value = 80.34
$ws.Range("E5").Value = 72.31
$ws.Range("C6").Value = 77
$ws.Range("E6").Value = 65.45

# --- Delete rows 10-13 (the old id = 3 group) first, so row indices for 7-9 stay put ---
$ws.Range("A10:F13").EntireRow.Delete()

# --- Rows 7-9 now become the id = 3 group with new project data ---
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "LLMGuard"
$ws.Range("C7").Value = 109.62
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 100
$ws.Range("F7").Value = 3

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Multi Model Data Analysis for Annotation of Human Activities"
$ws.Range("C8").Value = 102.85
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 100
$ws.Range("F8").Value = 3

$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "Website for the Literature Society of the college"
$ws.Range("C9").Value = 80.34
$ws.Range("D9").Value = 0.9
$ws.Range("E9").Value = 72.31
$ws.Range("F9").Value = 3
